# Data-driven test workbook: rename sheets and add an "error message" column
# to both the Login and Register data tables (register page gets a brand new
# column D with one error-message value per scenario row).

$wb = $excel.ActiveWorkbook

# --- Rename sheets -----------------------------------------------------
$wsLogin = $wb.Worksheets.Item(1)
$wsLogin.Name = "Login"

$wsRegister = $wb.Worksheets.Item(2)
$wsRegister.Name = "Register"

# --- Login sheet: add "error message" column (C) -----------------------
$wsLogin.Range("B1").Copy()
$wsLogin.Range("C1").PasteSpecial(-4122)
$wsLogin.Range("C1").Value = "error message"

$wsLogin.Range("B2").Copy()
$wsLogin.Range("C2").PasteSpecial(-4122)

$wsLogin.Range("B3").Copy()
$wsLogin.Range("C3").PasteSpecial(-4122)

$wsLogin.Range("B4").Copy()
$wsLogin.Range("C4").PasteSpecial(-4122)
$wsLogin.Range("C4").Value = "Invalid Username and Password"

# --- Register sheet: add "error message" column (D) --------------------
$wsRegister.Range("C1").Copy()
$wsRegister.Range("D1").PasteSpecial(-4122)
$wsRegister.Range("D1").Value = "error message"

$wsRegister.Range("C2").Copy()
$wsRegister.Range("D2").PasteSpecial(-4122)

$wsRegister.Range("C3").Copy()
$wsRegister.Range("D3").PasteSpecial(-4122)

$wsRegister.Range("C4").Copy()
$wsRegister.Range("D4").PasteSpecial(-4122)

$wsRegister.Range("C5").Copy()
$wsRegister.Range("D5").PasteSpecial(-4122)
$wsRegister.Range("D5").Value = "Your password can't be too similar to your other personal information."

$wsRegister.Range("C6").Copy()
$wsRegister.Range("D6").PasteSpecial(-4122)
$wsRegister.Range("D6").Value = "Your password must contain at least 8 characters."

$wsRegister.Range("C7").Copy()
$wsRegister.Range("D7").PasteSpecial(-4122)
$wsRegister.Range("D7").Value = "Your password can't be entirely numeric."

$wsRegister.Range("C8").Copy()
$wsRegister.Range("D8").PasteSpecial(-4122)
$wsRegister.Range("D8").Value = "password_mismatch:The two password fields didn't match."

$wsRegister.Range("C9").Copy()
$wsRegister.Range("D9").PasteSpecial(-4122)
$wsRegister.Range("D9").Value = "User already exists"
